# Insert a new row at position 11, shifting existing rows 11-64 down to 12-65,
# then populate the newly inserted row 11 with its data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 11 (shifts rows 11..64 -> 12..65)
$ws.Rows.Item(11).Insert()

# Fill in the data for the newly inserted row 11
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44883
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112027
$ws.Range("G11").Value = "Melón"
$ws.Range("H11").Value = "Tuna"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 28000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 29333
$ws.Range("N11").Value = "`$/caja 18 unidades"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 1630
$ws.Range("Q11").Value = 18
$ws.Range("R11").Value = "Hortaliza"
